$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure text cells keep their literal string representation (avoid Excel
# auto-converting numeric- or percent-looking strings into Number values).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.59%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.08%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.131'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.46%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07694'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.46%'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.306'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.87%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.622'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.21%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9226'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.73%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.458'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.03%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1217'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '20.72%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1825'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.12%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09150'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.47%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04323'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.11%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1050'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.51%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001246'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.09%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005846'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.10%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.346'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.18%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.941'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '5.40%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.42%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2678'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.86%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04052'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.97%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001263'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.99%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004097'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.99%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001268'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.61%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '24.61%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02473'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4.30%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05274'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.61%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007820'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.80%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1314'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.44%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006797'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.96%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001841'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.16%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008167'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.27%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3106'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.69%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006753'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.25%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.20%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2055'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '2,090.17%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004096'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-7.01%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.20%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001998'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.20%'
